$d = $word.ActiveDocument
$sec = $d.Sections(1)

# Footer 1: Pearson Edexcel logo (docPr id="3") -> rename image1.png to image2.png
$f1 = $sec.Footers(1)
$f1r = $f1.Range
for ($i = 1; $i -le $f1r.InlineShapes.Count; $i++) {
    $shp = $f1r.InlineShapes($i)
    if ($shp.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
        $shp.Name = "image2.png"
    }
}

# Footer 2: Pearson Edexcel logo (docPr id="2") -> rename image1.png to image2.png
$f2 = $sec.Footers(2)
$f2r = $f2.Range
for ($i = 1; $i -le $f2r.InlineShapes.Count; $i++) {
    $shp = $f2r.InlineShapes($i)
    if ($shp.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
        $shp.Name = "image2.png"
    }
}

# Header with BTec logo (docPr id="1") -> rename image2.jpg to image1.jpg
for ($hi = 1; $hi -le $sec.Headers.Count; $hi++) {
    $h = $sec.Headers($hi)
    if ($h.Exists) {
        $hr = $h.Range
        for ($i = 1; $i -le $hr.InlineShapes.Count; $i++) {
            $shp = $hr.InlineShapes($i)
            if ($shp.AlternativeText -eq "BTec_Logo-Orange") {
                $shp.Name = "image1.jpg"
            }
        }
    }
}
